$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.412.94'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -1.23%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.220.06'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.67%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.27'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.79%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '182.52'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.36%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.219.66'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.63%  '
$ws.Range("E10").Value = '  -3.58%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.52'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -3.03%  '
$ws.Range("E12").Value = '  -1.54%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.778.22'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.64%  '
$ws.Range("E14").Value = '  -0.05%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.68'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -3.63%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.474.76'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.14%  '
$ws.Range("E17").Value = '  -2.48%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.222.87'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.52%  '
$ws.Range("E19").Value = '  -2.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.39'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.85%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '393.89'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.54'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.52%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.01'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.64%  '
$ws.Range("E25").Value = '  -0.73%  '
$ws.Range("E26").Value = '  -3.56%  '
$ws.Range("E27").Value = '  -0.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.54'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.56%  '
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("E30").Value = '  -2.65%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.55'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -4.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '22.56'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.77%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.93'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -4.57%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("E35").Value = '  -3.01%  '
$ws.Range("E36").Value = '  -1.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.47'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -5.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.87'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.16%  '
$ws.Range("E39").Value = '  -4.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '26.15'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.50%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.55'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.67%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.49'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -4.40%  '
$ws.Range("E43").Value = '  -6.06%  '
$ws.Range("E44").Value = '  -1.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '40.49'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.587.38'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.57%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.47'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -4.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '332.88'
$ws.Range("D48").ClearFormats()
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0277'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.97%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.27'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.01%  '
$ws.Range("E51").Value = '  -2.09%  '
